$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Grille")

# --- Merge J42:J43 first so setting J42 below lands in the merged area ---
$ws.Range("J42:J43").Merge()

# --- Set cell values in the same order the original author entered them, so
#     shared-string indices line up with the source workbook. ---
$ws.Range("I2").Value = "Notes groupes"
$ws.Range("J2").Value = "Notes prof"
$ws.Range("J7").Value = "Voulait clicable au lieu de bouton"
$ws.Range("J12").Value = "En cas de DB down, pas de site"
$ws.Range("J19").Value = "hyper résumé"
$ws.Range("J21").Value = "ok dans excel donc bonus"
$ws.Range("J24").Value = "page dans titre, api dans task"
$ws.Range("J26").Value = "res.render sur le résultat, traite pas l'erreur"
$ws.Range("J30").Value = "prof va vérifier"
$ws.Range("J33").Value = "branches par features`ngael moins de commit du a sa manière de travailler`nfin de certaines branches encore a merger"
$ws.Range("J38").Value = "dans excel`nrésumé dans wiki"
$ws.Range("J39").Value = "non"
$ws.Range("J42").Value = "node.js bien : js // express à découvrir`nnode à part ressemble pas à js … impresion de faire du java … // pas clair`nSUPER ABSTRAIT // pour token require // communiquer backend front end se renseigner et ont mis en place des …tableaux // entre différentes pages on s'envoient des données… "
$ws.Range("J44").Value = "base générée en webstorm"
$ws.Range("J45").Value = "pas été tres en détails, semblent bien connaitre"
$ws.Range("J46").Value = "OK dans wiki directement"
$ws.Range("J47").Value = "format de réponse serait mieux du json directement"
$ws.Range("J48").Value = "manque login pas encore complet"
$ws.Range("J50").Value = "postman"
$ws.Range("J52").Value = "mysql ? Pq pas mariadb?"
$ws.Range("J53").Value = "Opensource pas clair"
$ws.Range("J54").Value = "comprennent leur connexion mais bof"
$ws.Range("J55").Value = "diagramme ER ok"
$ws.Range("J56").Value = "pas encore en ligne"
$ws.Range("J87").Value = "Attention en root et sans password"

# --- Formatting: wrap text for the long multi-line note cells ---
$ws.Range("J42:J43").WrapText = $true
$ws.Range("J42:J43").HorizontalAlignment = -4131
$ws.Range("J33").WrapText = $true
$ws.Range("J38").WrapText = $true

# --- Column J width ---
$ws.Columns.Item(10).ColumnWidth = 48.8

# --- Restore selection to match author final view ---
$ws.Activate()
$ws.Range("J88").Select()
